$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New values for the "costo_hh_us" (MB) column D, rows 2-28 (row 4 stays blank)
$newValues = @{
    2  = 27
    3  = 27
    5  = 25
    6  = 26
    7  = 26
    8  = 25
    9  = 26
    10 = 26
    11 = 28
    12 = 27
    13 = 27
    14 = 27
    15 = 25
    16 = 27
    17 = 26
    18 = 28
    19 = 26
    20 = 27
    21 = 27
    22 = 27
    23 = 27
    24 = 28
    25 = 27
    26 = 28
    27 = 28
    28 = 28
}

foreach ($row in $newValues.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $newValues[$row]
    $cell.ClearFormats()
    $cell.NumberFormat = "0"
    $cell.Font.Color = 0
}

# Row 4's D cell has no value but keeps the new number format
$blankCell = $ws.Cells.Item(4, 4)
$blankCell.ClearFormats()
$blankCell.NumberFormat = "0"
$blankCell.Font.Color = 0

# Update the active selection to match the saved workbook state
$ws.Range("F5").Select()
